$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Val)
    $Worksheet.Range($Address).NumberFormat = "@"
    $Worksheet.Range($Address).Value = $Val
    $Worksheet.Range($Address).Style = "Normal"
}

Set-TextValue $ws "D2" "42.964.80"
Set-TextValue $ws "E2" "  -1.19%  "
Set-TextValue $ws "D3" "2.311.20"
Set-TextValue $ws "E3" "  -0.25%  "
Set-TextValue $ws "D4" "0.998"
Set-TextValue $ws "E4" "  -0.24%  "
Set-TextValue $ws "D5" "301.73"
Set-TextValue $ws "E5" "  -2.38%  "
Set-TextValue $ws "D6" "99.24"
Set-TextValue $ws "E6" "  -6.39%  "
Set-TextValue $ws "D7" "0.504"
Set-TextValue $ws "E7" "  -4.20%  "
Set-TextValue $ws "D8" "0.999"
Set-TextValue $ws "E8" "  -0.16%  "
Set-TextValue $ws "D9" "0.502"
Set-TextValue $ws "E9" "  -3.48%  "
Set-TextValue $ws "D10" "34.78"
Set-TextValue $ws "E10" "  -4.21%  "
Set-TextValue $ws "D11" "0.0790"
Set-TextValue $ws "E11" "  -2.92%  "
Set-TextValue $ws "E12" "  +0.38%  "
Set-TextValue $ws "D13" "6.71"
Set-TextValue $ws "E13" "  -4.06%  "
Set-TextValue $ws "D14" "2.666.08"
Set-TextValue $ws "E14" "  -0.18%  "
Set-TextValue $ws "D15" "15.44"
Set-TextValue $ws "E15" "  +0.39%  "
Set-TextValue $ws "D16" "2.313.13"
Set-TextValue $ws "E16" "  -1.68%  "
Set-TextValue $ws "D17" "0.792"
Set-TextValue $ws "E17" "  -1.45%  "
Set-TextValue $ws "D18" "42.828.07"
Set-TextValue $ws "E18" "  -1.58%  "
Set-TextValue $ws "D19" "11.67"
Set-TextValue $ws "E19" "  -2.39%  "
Set-TextValue $ws "D20" "0.0₃0899"
Set-TextValue $ws "E20" "  -2.69%  "
Set-TextValue $ws "D21" "6.02"
Set-TextValue $ws "E21" "  -2.97%  "
Set-TextValue $ws "D22" "67.69"
Set-TextValue $ws "E22" "  -0.69%  "
Set-TextValue $ws "D23" "236.68"
Set-TextValue $ws "E23" "  -1.99%  "
Set-TextValue $ws "D24" "1.95"
Set-TextValue $ws "E24" "  -4.72%  "
Set-TextValue $ws "D25" "2.50"
Set-TextValue $ws "E25" "  -4.17%  "
Set-TextValue $ws "E26" "  +0.07%  "
Set-TextValue $ws "D27" "24.64"
Set-TextValue $ws "E27" "  -1.51%  "
Set-TextValue $ws "E28" "  -1.94%  "
Set-TextValue $ws "D29" "34.24"
Set-TextValue $ws "E29" "  -6.45%  "
Set-TextValue $ws "D30" "164.11"
Set-TextValue $ws "E30" "  +0.74%  "
Set-TextValue $ws "D31" "9.11"
Set-TextValue $ws "E31" "  -5.08%  "
Set-TextValue $ws "E32" "  -0.19%  "
Set-TextValue $ws "D33" "5.01"
Set-TextValue $ws "E33" "  -4.60%  "
Set-TextValue $ws "E34" "  -4.72%  "
Set-TextValue $ws "D35" "4.46"
Set-TextValue $ws "E35" "  -3.21%  "
Set-TextValue $ws "B36" "Celestia"
Set-TextValue $ws "C36" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws "D36" "16.65"
Set-TextValue $ws "E36" "  -9.18%  "
Set-TextValue $ws "B37" "Hedera"
Set-TextValue $ws "C37" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D37" "0.0695"
Set-TextValue $ws "E37" "  -5.67%  "
Set-TextValue $ws "D38" "2.89"
Set-TextValue $ws "E38" "  -4.39%  "
Set-TextValue $ws "D39" "1.79"
Set-TextValue $ws "E39" "  -4.15%  "
Set-TextValue $ws "D40" "0.100"
Set-TextValue $ws "E40" "  -5.75%  "
Set-TextValue $ws "E41" "  -4.21%  "
Set-TextValue $ws "D42" "2.48"
Set-TextValue $ws "E42" "  -0.14%  "
Set-TextValue $ws "D43" "1.960.91"
Set-TextValue $ws "E43" "  -0.24%  "
Set-TextValue $ws "D44" "0.0280"
Set-TextValue $ws "E44" "  -3.56%  "
Set-TextValue $ws "D45" "18.31"
Set-TextValue $ws "E45" "  -2.59%  "
Set-TextValue $ws "D46" "10.16"
Set-TextValue $ws "E46" "  -1.47%  "
Set-TextValue $ws "D47" "2.87"
Set-TextValue $ws "E47" "  -6.47%  "
Set-TextValue $ws "D48" "54.99"
Set-TextValue $ws "E48" "  -6.01%  "
Set-TextValue $ws "D49" "2.528.34"
Set-TextValue $ws "E49" "  -0.53%  "
Set-TextValue $ws "B50" "THORChain"
Set-TextValue $ws "C50" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws "D50" "4.68"
Set-TextValue $ws "E50" "  -1.98%  "
Set-TextValue $ws "B51" "HuobiToken"
Set-TextValue $ws "C51" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws "D51" "2.81"
Set-TextValue $ws "E51" "  -5.02%  "
